# "break out stock.yaml completed"
#
# 1) Append 18 new "day" rows (774-791) scraped on 22/10/2024 11:36:25.
# 2) Normalize the "week" sheet's bsecode column (D395:D413) from text to
#    numbers, matching the numeric-typed bsecode cells used everywhere else.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "day" sheet: append rows 774-791
# ---------------------------------------------------------------------
$day = $wb.Worksheets.Item("day")

$newRows = @(
    @(1, "OFSS", "Oracle Financial Services Software Limited", "532466", -2.73, 10890.7, 156837, "day", "22/10/2024 11:36:25"),
    @(2, "INDIGO", "Interglobe Aviation Limited", "539448", -1.45, 4524.4, 589153, "day", "22/10/2024 11:36:25"),
    @(3, "HDFCAMC", "HDFC Asset Management Company Ltd", "541729", -2, 4510.85, 567042, "day", "22/10/2024 11:36:25"),
    @(4, "BALKRISIND", "Balkrishna Industries Limited", "502355", -0.99, 2940.05, 152063, "day", "22/10/2024 11:36:25"),
    @(5, "OBEROIRLTY", "Oberoi Realty Limited", "533273", -2.26, 1950.45, 1764458, "day", "22/10/2024 11:36:25"),
    @(6, "MUTHOOTFIN", "Muthoot Finance Limited", "533398", -1.77, 1910.5, 357802, "day", "22/10/2024 11:36:25"),
    @(7, "VOLTAS", "Voltas Limited", "500575", -1.06, 1794.65, 984280, "day", "22/10/2024 11:36:25"),
    @(8, "COROMANDEL", "Coromandel International Limited", "506395", -2.2, 1569.9, 190581, "day", "22/10/2024 11:36:25"),
    @(9, "MFSL", "Max Financial Services Limited", "500271", -1.91, 1170.45, 1248285, "day", "22/10/2024 11:36:25"),
    @(10, "BSOFT", "Birlasoft Ltd", "532400", -3.47, 576.65, 2254673, "day", "22/10/2024 11:36:25"),
    @(11, "WIPRO", "Wipro Limited", "507685", -0.48, 545.45, 9239150, "day", "22/10/2024 11:36:25"),
    @(12, "RECLTD", "Rural Electrification Corporation Limited", "532955", -4.59, 508.05, 10871784, "day", "22/10/2024 11:36:25"),
    @(13, "APOLLOTYRE", "Apollo Tyres Limited", "500877", -2.59, 497.1, 1449938, "day", "22/10/2024 11:36:25"),
    @(14, "CHAMBLFERT", "Chambal Fertilizers & Chemicals Limited", "500085", -1.22, 481.3, 1396272, "day", "22/10/2024 11:36:25"),
    @(15, "PFC", "Power Finance Corporation Limited", "532810", -4.64, 442.4, 10178704, "day", "22/10/2024 11:36:25"),
    @(16, "HINDPETRO", "Hindustan Petroleum Corporation Limited", "500104", -3.14, 402.1, 5951077, "day", "22/10/2024 11:36:25"),
    @(17, "NMDC", "Nmdc Limited", "526371", -4.51, 215.41, 6708806, "day", "22/10/2024 11:36:25"),
    @(18, "FEDERALBNK", "The Federal Bank  Limited", "500469", -2.09, 189.34, 5731833, "day", "22/10/2024 11:36:25")
)

$startRow = 774
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $day.Cells.Item($r, 1).Value = $data[0]   # sr (number)
    $day.Cells.Item($r, 2).Value = $data[1]   # nsecode (text)
    $day.Cells.Item($r, 3).Value = $data[2]   # name (text)

    # bsecode is stored as TEXT for these freshly appended rows (even
    # though it is numeric-looking) - force text via NumberFormat so the
    # COM layer doesn't silently coerce the numeric-looking string.
    $bsecodeCell = $day.Cells.Item($r, 4)
    $bsecodeCell.NumberFormat = "@"
    $bsecodeCell.Value = $data[3]

    $day.Cells.Item($r, 5).Value = $data[4]   # per_chg (number)
    $day.Cells.Item($r, 6).Value = $data[5]   # close (number)
    $day.Cells.Item($r, 7).Value = $data[6]   # volume (number)
    $day.Cells.Item($r, 8).Value = $data[7]   # timeframe (text)
    $day.Cells.Item($r, 9).Value = $data[8]   # Date Time (text)
}

# ---------------------------------------------------------------------
# 2) "week" sheet: D395:D413 bsecode text -> number
# ---------------------------------------------------------------------
$week = $wb.Worksheets.Item("week")

$bsecodes = @(500550, 539523, 539448, 532644, 500495, 542726, 533398, 524804, 532424, 542830, 511196, 500112, 543066, 540222, 532555, 535755, 526371, 500469, 532483)

$row = 395
for ($i = 0; $i -lt $bsecodes.Length; $i++) {
    $week.Cells.Item($row, 4).Value = $bsecodes[$i]
    $row++
}

Write-Host "new_ph_pl: appended day rows 774-791, normalized week!D395:D413 to numeric"
